# Final participants update
# - fix Age for participant 62 (row 63)
# - fill in Gender/Age for the last 4 participants (rows 70-73)
# - add a new "Notes" column (D) with per-participant notes
# - format the header row (bold, centered, filled) and widen the Notes column
# - adjust sheet selection / view

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- fix a data entry mistake -------------------------------------------------
$ws.Range("C63").Value = 19

# --- fill in the four participants that were missing Gender/Age --------------
$ws.Range("B70").Value = "Female"
$ws.Range("C70").Value = 19

$ws.Range("B71").Value = "Female"
$ws.Range("C71").Value = 20

$ws.Range("B72").Value = "Male"
$ws.Range("C72").Value = 20

$ws.Range("B73").Value = "Female"
$ws.Range("C73").Value = 20

# --- add the Notes column ------------------------------------------------------
$ws.Range("D1").Value = "Notes"

$ws.Range("D3").Value = "glasses"
$ws.Range("D13").Value = "Contact lenses"
$ws.Range("D17").Value = "glasses"
$ws.Range("D28").Value = "glasses; left eye tracked"
$ws.Range("D31").Value = "glasses"
$ws.Range("D32").Value = "glasses"
$ws.Range("D34").Value = "left eye tracked"
$ws.Range("D35").Value = "glasses; left eye tracked"
$ws.Range("D39").Value = "left eye tracked"
$ws.Range("D61").Value = "Contact lenses"
$ws.Range("D63").Value = "glasses"
$ws.Range("D69").Value = "glasses"

# --- header formatting ----------------------------------------------------------
$header = $ws.Range("A1:D1")
$header.Font.Bold = $true
$header.Interior.ThemeColor = 5
$header.Interior.TintAndShade = 0.59999389629810485
$header.HorizontalAlignment = -4108

# --- column widths ---------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 23.6

# --- selection / view -------------------------------------------------------------
$ws.Range("I9").Select()
